$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.942.48'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '2.907.77'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '569.93'
$ws.Range('E5').Value = '  -3.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.23'
$ws.Range('E6').Value = '  -1.52%  '
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D8').Value = '2.905.89'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('B9').Value = 'XRP'
$ws.Range('C9').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.501'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.99'
$ws.Range('E10').Value = '  -1.00%  '
$ws.Range('E11').Value = '  -2.46%  '
$ws.Range('E12').Value = '  -0.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000231'
$ws.Range('E13').Value = '  -1.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.32'
$ws.Range('E14').Value = '  +0.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.126'
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('D16').Value = '3.388.61'
$ws.Range('E16').Value = '  -0.71%  '
$ws.Range('D17').Value = '61.883.00'
$ws.Range('E17').Value = '  -0.99%  '
$ws.Range('D18').Value = '2.904.78'
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.50'
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '429.88'
$ws.Range('E20').Value = '  -0.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.98'
$ws.Range('E21').Value = '  -2.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.651'
$ws.Range('E22').Value = '  -0.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.87'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.89'
$ws.Range('E24').Value = '  -1.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.03'
$ws.Range('E25').Value = '  +2.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.20'
$ws.Range('E26').Value = '  -7.28%  '
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.03'
$ws.Range('E28').Value = '  -2.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000112'
$ws.Range('E29').Value = '  +10.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.99'
$ws.Range('E30').Value = '  -2.89%  '
$ws.Range('E31').Value = '  -2.44%  '
$ws.Range('E32').Value = '  -4.93%  '
$ws.Range('E33').Value = '  -0.21%  '
$ws.Range('E34').Value = '  -1.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '25.62'
$ws.Range('E35').Value = '  -1.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.955'
$ws.Range('E36').Value = '  -1.96%  '
$ws.Range('E37').Value = '  -1.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '48.82'
$ws.Range('E38').Value = '  -0.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.89'
$ws.Range('E39').Value = '  -5.65%  '
$ws.Range('E40').Value = '  -4.90%  '
$ws.Range('E41').Value = '  -0.58%  '
$ws.Range('B42').Value = 'Arweave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '41.19'
$ws.Range('E42').Value = '  +4.35%  '
$ws.Range('B43').Value = 'Cosmos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.16'
$ws.Range('E43').Value = '  -1.94%  '
$ws.Range('E44').Value = '  -2.93%  '
$ws.Range('D45').Value = '2.703.72'
$ws.Range('E45').Value = '  +0.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0336'
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '131.93'
$ws.Range('E47').Value = '  -2.91%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '347.08'
$ws.Range('E48').Value = '  -1.59%  '
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.58'
$ws.Range('E51').Value = '  -3.53%  '
